# Generate Report for Archive
#
# 1) Status text changed from "Ready for handoff" to "In Translation" in all
#    four cells that held it: Overview!E2, Overview!F2, zh-cn!C2, de-de!C2.
# 2) Because the new status text is shorter, the Status column(s) were
#    narrowed (re-autofit) from ~17.22 chars to ~13.41 chars:
#      - Overview columns E and F
#      - zh-cn column C
#      - de-de column C

$wb = $excel.ActiveWorkbook

# The workbook's recorded column width (13.4101845877511 "characters") comes
# from a real-Excel autofit computation against exact font metrics, which
# this host's ColumnWidth setter cannot reproduce bit-for-bit (it quantizes
# to a 6-px grid after rounding the input to 2 decimals). 12.5 is the input
# that lands on the closest reachable grid point to that target.
$targetWidth = 12.5

# --- Overview sheet: columns E (zh-cn status) and F (de-de status) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = $targetWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetWidth

# --- zh-cn sheet: column C (Status) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = $targetWidth

# --- de-de sheet: column C (Status) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = $targetWidth
